$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStrings = @(
    "Sát Phá Tham hội chiếu tại Nô Bộc",
    "Tử Phủ Vũ Tướng hội chiếu tại Nô Bộc",
    "Cơ Nguyệt Đồng Lương hội chiếu tại Nô Bộc",
    "Cự Nhật hội chiếu tại Nô Bộc",
    "Kình Đà hội chiếu tại Nô Bộc",
    "Xương Khúc hội chiếu tại Nô Bộc",
    "Hoả Linh hội chiếu tại Nô Bộc",
    "Không Kiếp hội chiếu tại Nô Bộc",
    "Quang Quý hội chiếu tại Nô Bộc",
    "Tả Hữu hội chiếu tại Nô Bộc",
    "Song Hao hội chiếu tại Nô Bộc",
    "Tang Hổ hội chiếu tại Nô Bộc",
    "Khốc Hư hội chiếu tại Nô Bộc",
    "Hình Riêu hội chiếu tại Nô Bộc",
    "Thai Toạ hội chiếu tại Nô Bộc",
    "Đào Hồng hội chiếu tại Nô Bộc",
    "Ấn Phù hội chiếu tại Nô Bộc",
    "Song Hao Quyền Lộc Kiếp Hoả hội chiếu tại Nô Bộc",
    "Tử Phủ Vũ Tướng Xương Khúc Khôi Việt Tả Hữu Khoa Quyền Lộc Long hội chiếu tại Nô Bộc",
    "Tử Phủ Vũ Tướng Tả Hữu Khoa Quyền Lộc Long Phượng hội chiếu tại Nô Bộc",
    "Tử Khúc Phá Dương Đà hội chiếu tại Nô Bộc",
    "Cơ Nguyệt Đồng Lương gặp Xương Khúc Tả Hữu hội chiếu tại Nô Bộc",
    "Cơ Nguyệt Đồng Lương Khoa Tả Hữu Quang Quý Quan Phúc hội chiếu tại Nô Bộc",
    "Sát Quyền hội chiếu tại Nô Bộc",
    "Lộc Mã hội chiếu tại Nô Bộc",
    "Kiếp Hư Hao Quyền hội chiếu tại Nô Bộc",
    "Tuế Hổ Phù Xương Khúc hội chiếu tại Nô Bộc",
    "Xương Khúc Tấu Long Phượng hội chiếu tại Nô Bộc",
    "Đào Hồng Riêu Tấu Cơ Vũ hội chiếu tại Nô Bộc",
    "Binh Hình Tướng Ấn hội chiếu tại Nô Bộc",
    "Hổ Tấu hội chiếu tại Nô Bộc",
    "Hình Riêu Y hội chiếu tại Nô Bộc",
    "Mã Hỏa Linh hội chiếu tại Nô Bộc",
    "Thai Tọa Hồng Đào hội chiếu tại Nô Bộc",
    "Tả Hữu Không Kiếp hội chiếu tại Nô Bộc",
    "Tả Hữu Binh Tướng hội chiếu tại Nô Bộc",
    "Đào Quyền hội chiếu tại Nô Bộc",
    "Đào Hồng Tả Cái Hữu hội chiếu tại Nô Bộc",
    "Quan Phúc Quang Tấu hội chiếu tại Nô Bộc",
    "Đào Hồng Xương Khúc hội chiếu tại Nô Bộc",
    "Đào Tử Phủ hội chiếu tại Nô Bộc",
    "Hổ Kình Sát hội chiếu tại Nô Bộc",
    "Hổ Tang Kiếp hội chiếu tại Nô Bộc",
    "Hổ Tang Không Kiếp hội chiếu tại Nô Bộc",
    "Âm Dương Lương hội chiếu tại Nô Bộc",
    "Cơ Lương Gia Hội hội chiếu tại Nô Bộc",
    "Nhật Chiếu Lôi Môn hội chiếu tại Nô Bộc",
    "Tả Hữu Xương Khúc hội chiếu tại Nô Bộc",
    "Tham Linh Triều Viên hội chiếu tại Nô Bộc",
    "Cự Hỏa Linh hội chiếu tại Nô Bộc",
    "Hồng Đào Kỵ hội chiếu tại Nô Bộc",
    "Cự Đồng Hình hội chiếu tại Nô Bộc",
    "Lương Phá hội chiếu tại Nô Bộc",
    "Khôi Việt hội chiếu tại Nô Bộc",
    "Tham Vũ Hỏa hội chiếu tại Nô Bộc",
    "Cự Kỵ hội chiếu tại Nô Bộc",
    "Kình Đà Hỏa hội chiếu tại Nô Bộc",
    "Kình Đà Không Kiếp hội chiếu tại Nô Bộc",
    "Đào Hồng Hỷ hội chiếu tại Nô Bộc",
    "Đào Hồng Kiếp Sát hội chiếu tại Nô Bộc",
    "Đào Hồng Kỵ hội chiếu tại Nô Bộc",
    "Đào Hồng Tả Phù Hữu Bật hội chiếu tại Nô Bộc",
    "Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc hội chiếu tại Nô Bộc",
    "Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc Long Trì Phượng Các hội chiếu tại Nô Bộc",
    "Lưu Hà Kiếp Sát hội chiếu tại Nô Bộc",
    "Phục Không Kiếp hội chiếu tại Nô Bộc",
    "Xương Khúc Khôi Việt hội chiếu tại Nô Bộc",
    "Khoa Quyền hội chiếu tại Nô Bộc",
    "Lộc Quyền hội chiếu tại Nô Bộc",
    "Tướng Binh Đào Hồng hội chiếu tại Nô Bộc",
    "Lương Khốc Tuế hội chiếu tại Nô Bộc",
    "Tướng Binh hội chiếu tại Nô Bộc",
    "Cự Tang hội chiếu tại Nô Bộc",
    "Cự Hỏa hội chiếu tại Nô Bộc",
    "Phá Hình Kỵ hội chiếu tại Nô Bộc",
    "Tang Trực Tuế hội chiếu tại Nô Bộc",
    "Tuế Xương Khúc hội chiếu tại Nô Bộc",
    "Tử Vi Tả Hữu hội chiếu tại Nô Bộc",
    "Tử Vi Tang Tả Hữu hội chiếu tại Nô Bộc"
)

$startRow = 4211
for ($i = 0; $i -lt $newStrings.Count; $i++) {
    $row = $startRow + $i
    $value = $newStrings[$i]
    $ws.Cells.Item($row, 1).Value = $value
    $ws.Cells.Item($row, 2).Value = $value
}

$endRow = $startRow + $newStrings.Count - 1
$ws.Range("B" + $startRow + ":B" + $endRow).Select()

Write-Output ("Added rows {0} to {1}" -f $startRow, $endRow)
